# Applies the diff:
#  - Row 1 "99.97" -> "0M"
#  - Row 2 "0"     -> "0M"
#  - Row 3 "15"    -> "0M"
#  - Insert 10 new one-line rows after row 3 with values:
#       101, 0.00003, 0.00022, 0.00006, 0.00002, 0.00007, 0.00007, 0.00007, 0.00479, 100.0
#  - The row that used to hold "100<tab>0.00003<tab>...<tab>100.0" now holds just "99.97"
#  - The row that used to hold "1<tab>0.00007<tab>...<tab>100.0" now holds just "0"
#  - The last (previously empty) row now holds "15"

$d = $word.ActiveDocument
$t = $d.Tables(1)

function Set-CellText($cell, [string]$newText) {
    $r = $cell.Range
    # Exclude the trailing cell-mark (and any paragraph mark) from the
    # range so assigning .Text replaces the content without clobbering
    # the cell/row structure.
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $newText
}

# 1) First three rows become "0M"
Set-CellText $t.Rows(1).Cells(1) "0M"
Set-CellText $t.Rows(2).Cells(1) "0M"
Set-CellText $t.Rows(3).Cells(1) "0M"

# 2) Insert 10 new rows right after row 3, one value per row.
#    Rows.Add(BeforeRow) always inserts immediately above the same
#    anchor row, so build them in reverse order to end up with the
#    values in the intended top-to-bottom order.
$newValues = @("101", "0.00003", "0.00022", "0.00006", "0.00002", "0.00007", "0.00007", "0.00007", "0.00479", "100.0")
$insertBeforeRow = $t.Rows(4)
for ($i = $newValues.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($insertBeforeRow)
    Set-CellText $newRow.Cells(1) $newValues[$i]
}

# 3) The two "wide" rows (originally tab-separated run lists) collapse to
#    a single value each; they are now rows 44 and 45 after the 10-row
#    insertion above (34 + 10, 35 + 10).
Set-CellText $t.Rows(44).Cells(1) "99.97"
Set-CellText $t.Rows(45).Cells(1) "0"

# 4) The final row (previously an empty run) gets the text "15"
Set-CellText $t.Rows(46).Cells(1) "15"
